# Update "想去人数" (F column) figures for both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 78
$ws1.Range("F5").Value = 3084
$ws1.Range("F7").Value = 2515
$ws1.Range("F10").Value = 4
$ws1.Range("F11").Value = 1275
$ws1.Range("F13").Value = 53
$ws1.Range("F15").Value = 1135
$ws1.Range("F16").Value = 323
$ws1.Range("F17").Value = 319
$ws1.Range("F19").Value = 26
$ws1.Range("F23").Value = 1094
$ws1.Range("F24").Value = 15
$ws1.Range("F25").Value = 259

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 78
$ws4.Range("F5").Value = 3084
$ws4.Range("F7").Value = 2515
$ws4.Range("F10").Value = 4
$ws4.Range("F11").Value = 1275
$ws4.Range("F13").Value = 53
$ws4.Range("F15").Value = 1135
$ws4.Range("F16").Value = 323
$ws4.Range("F17").Value = 319
$ws4.Range("F19").Value = 26
$ws4.Range("F23").Value = 1106
$ws4.Range("F24").Value = 15
$ws4.Range("F25").Value = 259
